$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.738.56'
$ws.Range('E2').Value = '  +1.25%  '
$ws.Range('D3').Value = '3.740.79'
$ws.Range('E3').Value = '  -1.58%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.12%  '
$ws.Range('D7').Value = '3.736.42'
$ws.Range('E7').Value = '  -1.66%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +0.70%  '
$ws.Range('E10').Value = '  +2.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.34'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.462'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.19'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('D15').Value = '4.363.58'
$ws.Range('E15').Value = '  -1.68%  '
$ws.Range('D16').Value = '3.730.88'
$ws.Range('E16').Value = '  -1.81%  '
$ws.Range('D17').Value = '68.743.71'
$ws.Range('E17').Value = '  +1.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.30'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.76'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +17.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '494.27'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.728'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('E26').Value = '  -3.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.43'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.54%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.57'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.98'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('E32').Value = '  +0.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.70'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.78%  '
$ws.Range('D34').Value = '3.882.92'
$ws.Range('E34').Value = '  -1.57%  '
$ws.Range('D35').Value = '3.672.01'
$ws.Range('E35').Value = '  -1.73%  '
$ws.Range('E36').Value = '  -1.42%  '
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.83'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('E40').Value = '  +0.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.326'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '438.29'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '48.93'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.21%  '
$ws.Range('E44').Value = '  -1.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.87'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.51'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.81%  '
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.74'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.16'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.70%  '
$ws.Range('D50').Value = '2.795.12'
$ws.Range('E50').Value = '  -1.81%  '
